$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.418.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.57%  "
$ws.Range("D3").Value = "'3.408.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.44%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'587.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").Value = "'181.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.86%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.45%  "
$ws.Range("E9").Value = "  +8.65%  "
$ws.Range("D10").Value = "'0.594"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.63%  "
$ws.Range("D11").Value = "'48.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.07%  "
$ws.Range("D12").Value = "'0.0000284"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.74%  "
$ws.Range("D13").Value = "'683.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("E14").Value = "  +3.65%  "
$ws.Range("D15").Value = "'3.949.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").Value = "'69.468.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.67%  "
$ws.Range("D17").Value = "'3.404.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.69%  "
$ws.Range("E18").Value = "  +1.60%  "
$ws.Range("D19").Value = "'17.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("D20").Value = "'11.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.05%  "
$ws.Range("D21").Value = "'0.910"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("D22").Value = "'17.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.49%  "
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("D24").Value = "'103.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.88%  "
$ws.Range("D25").Value = "'3.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("E26").Value = "  +2.45%  "
$ws.Range("E27").Value = "  +3.91%  "
$ws.Range("D28").Value = "'34.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.88%  "
$ws.Range("D29").Value = "'8.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.13%  "
$ws.Range("D30").Value = "'7.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("D31").Value = "'11.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.01%  "
$ws.Range("D32").Value = "'558.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.85%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.107"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.68%  "
$ws.Range("B34").Value = "dogwifhat"
$ws.Range("C34").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D34").Value = "'3.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.18%  "
$ws.Range("D35").Value = "'58.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.99%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "'3.673.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("D38").Value = "'0.141"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.02%  "
$ws.Range("D39").Value = "'35.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.44%  "
$ws.Range("D40").Value = "'0.0₃0725"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.74%  "
$ws.Range("E41").Value = "  +3.53%  "
$ws.Range("E42").Value = "  +3.37%  "
$ws.Range("D43").Value = "'0.341"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.11%  "
$ws.Range("E44").Value = "  +5.11%  "
$ws.Range("D45").Value = "'3.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("E46").Value = "  +2.32%  "
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("D48").Value = "'1.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.72%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").Value = "'134.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("D51").Value = "'2.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.07%  "
